# Update "想去人数" (people interested) counts on the "展览" and "全部类型" sheets
# to reflect the latest generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 10707
$ws1.Range("F5").Value = 698
$ws1.Range("F6").Value = 497

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 10707
$ws4.Range("F5").Value = 698
$ws4.Range("F7").Value = 497
